# Updating the staging env testdata
# Sheet "OldImportLogic" test values move from "_1" (TestAutomation_1) to
# "_2" (TestAutomation_2) in its Name/Name_radio_button/report columns.

$wb = $excel.ActiveWorkbook

$oldSheet = $wb.Worksheets.Item("OldImportLogic")

$oldSheet.Range("H2").Value = "StandardExcelReport-OldImportLogic_2-TestAutomation_2-Quality of Life-2023_"
$oldSheet.Range("H3").Value = "ExcelReport-OldImportLogic_2-TestAutomation_2-Quality of Life-"
$oldSheet.Range("H4").Value = "WordReport-OldImportLogic_2 - TestAutomation_2-Quality of Life-"
$oldSheet.Range("B2").Value = "OldImportLogic_2 - TestAutomation_2_radio_button"
$oldSheet.Range("A2").Value = "OldImportLogic_2 - TestAutomation_2"

# Make OldImportLogic the active sheet / tab, with A2 selected.
$oldSheet.Activate()
$oldSheet.Range("A2").Select()

# NewImportLogic selection moves to B2 (and the saved scroll position resets).
$newSheet = $wb.Worksheets.Item("NewImportLogic")
$newSheet.Range("B2").Select()

$oldSheet.Activate()
